$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting/row of the last existing entry (row 8: luigi/password/docente)
# and create a new row for "marco" with role "docente" (a user-created test)
$ws.Range("D8:F8").Copy()
$ws.Range("D9:F9").PasteSpecial(-4104)  # xlPasteAll

$ws.Range("D9").Value = "marco"
$ws.Range("E9").Value = "password"
$ws.Range("F9").Value = "docente"

$ws.Range("F19").Select()
